$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 10, pushing the existing rows 10-15 down to 11-16.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new review data.
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = ""
$ws.Cells.Item(10, 3).Value = 45919.53137916666
$ws.Cells.Item(10, 3).NumberFormat = $ws.Cells.Item(11, 3).NumberFormat
$ws.Cells.Item(10, 4).Value = "YWMxNDdmNmMtMThiZS00MDJmLThmYTEtM2E4MGZlYzY0MTRiOjU3MDE2"
